$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'68.506.97"
$ws.Cells.Item(2, 5).Value = "  +0.29%  "

$ws.Cells.Item(3, 4).Value = "'3.759.40"
$ws.Cells.Item(3, 5).Value = "  -0.93%  "

$ws.Cells.Item(4, 5).Value = "  +0.00%  "

$ws.Cells.Item(5, 4).Value = "'593.33"
$ws.Cells.Item(5, 5).Value = "  -0.81%  "

$ws.Cells.Item(6, 4).Value = "'167.07"
$ws.Cells.Item(6, 5).Value = "  -1.99%  "

$ws.Cells.Item(7, 4).Value = "'3.757.57"
$ws.Cells.Item(7, 5).Value = "  -0.98%  "

$ws.Cells.Item(8, 5).Value = "  -0.10%  "

$ws.Cells.Item(9, 5).Value = "  -1.14%  "

$ws.Cells.Item(10, 5).Value = "  -3.11%  "

$ws.Cells.Item(11, 5).Value = "  -1.68%  "

$ws.Cells.Item(12, 5).Value = "  -1.35%  "

$ws.Cells.Item(13, 5).Value = "  -7.50%  "

$ws.Cells.Item(14, 4).Value = "'36.13"
$ws.Cells.Item(14, 5).Value = "  -2.25%  "

$ws.Cells.Item(15, 4).Value = "'4.391.10"
$ws.Cells.Item(15, 5).Value = "  -0.96%  "

$ws.Cells.Item(16, 4).Value = "'3.752.56"
$ws.Cells.Item(16, 5).Value = "  -1.04%  "

$ws.Cells.Item(17, 4).Value = "'68.459.59"
$ws.Cells.Item(17, 5).Value = "  +0.20%  "

$ws.Cells.Item(18, 4).Value = "'18.01"
$ws.Cells.Item(18, 5).Value = "  -3.58%  "

$ws.Cells.Item(19, 5).Value = "  +0.74%  "

$ws.Cells.Item(20, 4).Value = "'6.99"
$ws.Cells.Item(20, 5).Value = "  -3.10%  "

$ws.Cells.Item(21, 5).Value = "  +1.00%  "

$ws.Cells.Item(22, 4).Value = "'465.07"
$ws.Cells.Item(22, 5).Value = "  -0.85%  "

$ws.Cells.Item(23, 5).Value = "  -3.23%  "

$ws.Cells.Item(24, 5).Value = "  -1.88%  "

$ws.Cells.Item(25, 4).Value = "'84.08"
$ws.Cells.Item(25, 5).Value = "  +0.19%  "

$ws.Cells.Item(26, 4).Value = "'2.18"
$ws.Cells.Item(26, 5).Value = "  -3.68%  "

$ws.Cells.Item(27, 4).Value = "'11.96"
$ws.Cells.Item(27, 5).Value = "  -2.01%  "

$ws.Cells.Item(28, 4).Value = "'10.05"
$ws.Cells.Item(28, 5).Value = "  -4.17%  "

$ws.Cells.Item(29, 5).Value = "  -0.13%  "

$ws.Cells.Item(30, 4).Value = "'3.907.79"
$ws.Cells.Item(30, 5).Value = "  -0.94%  "

$ws.Cells.Item(31, 5).Value = "  -4.97%  "

$ws.Cells.Item(32, 4).Value = "'7.34"
$ws.Cells.Item(32, 5).Value = "  -3.85%  "

$ws.Cells.Item(33, 5).Value = "  -2.02%  "

$ws.Cells.Item(34, 4).Value = "'2.18"
$ws.Cells.Item(34, 5).Value = "  -2.81%  "

$ws.Cells.Item(35, 5).Value = "  -1.13%  "

$ws.Cells.Item(36, 4).Value = "'0.998"

$ws.Cells.Item(37, 4).Value = "'3.717.31"
$ws.Cells.Item(37, 5).Value = "  -1.02%  "

$ws.Cells.Item(38, 4).Value = "'0.101"
$ws.Cells.Item(38, 5).Value = "  -3.79%  "

$ws.Cells.Item(39, 5).Value = "  -10.48%  "

$ws.Cells.Item(40, 2).Value = "Mantle"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(40, 4).Value = "'1.00"
$ws.Cells.Item(40, 5).Value = "  -1.10%  "

$ws.Cells.Item(41, 2).Value = "Kaspa"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(41, 4).Value = "'0.137"
$ws.Cells.Item(41, 5).Value = "  -1.77%  "

$ws.Cells.Item(42, 4).Value = "'5.79"
$ws.Cells.Item(42, 5).Value = "  -1.65%  "

$ws.Cells.Item(43, 5).Value = "  -0.12%  "

$ws.Cells.Item(45, 4).Value = "'43.97"
$ws.Cells.Item(45, 5).Value = "  +8.48%  "

$ws.Cells.Item(46, 4).Value = "'0.302"
$ws.Cells.Item(46, 5).Value = "  -3.92%  "

$ws.Cells.Item(47, 4).Value = "'46.75"
$ws.Cells.Item(47, 5).Value = "  +2.18%  "

$ws.Cells.Item(48, 5).Value = "  -2.21%  "

$ws.Cells.Item(49, 5).Value = "  -2.51%  "

$ws.Cells.Item(50, 4).Value = "'145.51"
$ws.Cells.Item(50, 5).Value = "  +1.66%  "

$ws.Cells.Item(51, 4).Value = "'389.42"
$ws.Cells.Item(51, 5).Value = "  -4.17%  "
